# Fix 15: Mapping Table 1 was missing a row for the "Urgent" source element
# that maps to MedicationRequest.priority with no source path on the
# duplicate/continuation row, plus refresh the Metadata "Date" value.

$wb = $excel.ActiveWorkbook

# --- 1. Update the metadata "Date" value on the Metadata sheet (B8) ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Cells.Item(8, 2).Value = "2025-01-11T17:10:04+00:00"

# --- 2. Mapping Table 1: insert a new row above row 5 ---
# Old row 5 ("Urgent" -> MedicationRequest.priority, a duplicate of row 4)
# becomes row 6, and both the newly inserted row 5 and the shifted row 6
# keep column D ("MedicationRequest.priority") but have column A (source
# path) cleared, matching the other grouped/continuation rows in the sheet.
$ws = $wb.Worksheets.Item("Mapping Table 1")

$ws.Rows.Item(5).Insert()
$ws.Range("A4:E4").Copy($ws.Range("A5:E5"))
$ws.Range("A5").ClearContents()
$ws.Range("A6").ClearContents()
